# Re-apply the results of re-running the stochastic loop that produced columns B and C.
# Column A and the header row (row 1) are untouched; only B2:C67 change value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$bValues = @(
    -0.2997831749105829,
    0.2213465968674888,
    0.1638914582948199,
    -0.04920836358425688,
    0.159229155010963,
    -0.4370824554597749,
    -0.2795541471843749,
    -0.4045438831751386,
    0.3305883237285109,
    -0.2009811357140011,
    -0.08134847341600106,
    0.002828604961007235,
    0.03936472424218167,
    -0.0784770619951347,
    0.2663150473889699,
    0.5817300000004855,
    0.0566309450292076,
    0.4030029950796304,
    0.24031889531072,
    0.4448595172459415,
    0.3508498213631931,
    -0.06541488236276555,
    4.532737832975839,
    0.5059606197432779,
    0.4388496050477586,
    0.3137994307768761,
    1.11590431386523,
    5.807625713805196,
    0.9916957704610752,
    -0.1604563312911147,
    0.8187632780202023,
    0.8990474586502597,
    -0.6292648004388127,
    0.8319206323793923,
    0.7856639771418961,
    0.7657448627624103,
    0.7594274471392801,
    0.5828628443902168,
    0.7545195054032509,
    0.5736053420197791,
    0.7225964188402608,
    0.7337029777463997,
    0.6742232400413685,
    0.6598447411032125,
    -1.250866761297211,
    -0.9684707031194607,
    -0.8614117227282849,
    -0.6303261647878988,
    -0.04745076767031046,
    -0.8505919969084399,
    -0.8505919969084399,
    -1.078203804713423,
    -0.1856291984996226,
    -0.9879768661455617,
    -0.8902972652459165,
    -0.9630652003590905,
    -1.164950966242617,
    -0.8680778599779868,
    -0.5206320563400149,
    0.3668523431147675,
    -1.252052505985923,
    -0.7626231060648025,
    -0.8999310999636885,
    -0.1389193221111201,
    -0.8312520780885234,
    -0.8328756808165422
)

$cValues = @(
    -0.2823340717132394,
    0.2315938719084131,
    0.1769144217514632,
    -0.03974547495026053,
    0.1749800295807246,
    -0.4265180793813819,
    -0.2573482178669795,
    -0.3937775289867715,
    0.3479642112326698,
    -0.1929657026794722,
    -0.05604147822376919,
    0.003674876464929709,
    0.07966766727242915,
    -0.02901894616415965,
    0.3327828794176409,
    0.6130156015452823,
    0.03819302181020207,
    0.4114595772555285,
    0.2998439411163769,
    0.5114454521818258,
    0.3871067794613055,
    -0.03357779522570183,
    4.549906949876535,
    0.4812838771670367,
    0.4153055876961433,
    0.2918728308909141,
    1.083150543661359,
    5.346451823179435,
    0.9482437885034567,
    -0.2160094535843192,
    0.781030889542029,
    0.8787370845010629,
    -0.650701299910941,
    0.8218666675194118,
    0.7693616306851152,
    0.7455412835311658,
    0.7362012578158906,
    0.5806585922460081,
    0.7531369528812434,
    0.5659644968844423,
    0.6926609629241252,
    0.717858937027485,
    0.6672250598973204,
    0.6428954243687939,
    -1.254761354839859,
    -0.9733229809055314,
    -0.8672692748034503,
    -0.6335504398257086,
    -0.04784744882433194,
    -0.8541424023608495,
    -0.8541424023608495,
    -1.092381861724099,
    -0.1838936610819277,
    -0.9915399417583468,
    -0.8842722965029342,
    -0.9459093158685727,
    -1.133883355160749,
    -0.8499787085037546,
    -0.49453003157713,
    0.3700907487945039,
    -1.229094600639659,
    -0.7263265570317615,
    -0.8911224901163964,
    -0.1133190903766008,
    -0.8004129880306827,
    -0.7906371324403927
)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}
